$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Puerro" (leek) price-history rows for Vega Central Mapocho de Santiago. ---
# Weekly refresh: three new sampling dates were folded into the existing
# history (rows 26, 40 and 50 below), pushing the remaining rows down and
# extending the sheet from 58 to 61 data rows.

$firstRow = 26
$lastRow  = 61

# Columns A,B,C,E,F,G,H,I,N,O,Q,R are identical on every data row in this
# sheet (same market/region/category), so they are (re)written uniformly.
$ws.Range("A$firstRow`:A$lastRow").Value = 9
$ws.Range("B$firstRow`:B$lastRow").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C$firstRow`:C$lastRow").Value = 'Metropolitana'
$ws.Range("E$firstRow`:E$lastRow").Value = 13
$ws.Range("F$firstRow`:F$lastRow").Value = 100112005
$ws.Range("G$firstRow`:G$lastRow").Value = 'Puerro'
$ws.Range("H$firstRow`:H$lastRow").Value = 'Sin especificar'
$ws.Range("I$firstRow`:I$lastRow").Value = 'Primera'
$ws.Range("N$firstRow`:N$lastRow").Value = '$/paquete 20 unidades'
$ws.Range("O$firstRow`:O$lastRow").Value = 'Provincia de Chacabuco'
$ws.Range("Q$firstRow`:Q$lastRow").Value = 20
$ws.Range("R$firstRow`:R$lastRow").Value = 'Hortaliza'

# Per-row varying data: Fecha(D), Volumen(J), Precio minimo(K), Precio
# maximo(L), Precio promedio ponderado(M), Precio $/Kg(P).
$data = @(
    @(44428, 97, 8000, 9000, 8505, 425),
    @(44162, 50, 8000, 8000, 8000, 400),
    @(44342, 160, 8000, 8000, 8000, 400),
    @(44384, 160, 8000, 9000, 8500, 425),
    @(44363, 160, 8000, 8000, 8000, 400),
    @(44349, 130, 8000, 8000, 8000, 400),
    @(44421, 180, 7000, 8000, 7500, 375),
    @(44273, 70, 8000, 8000, 8000, 400),
    @(44426, 97, 7000, 8000, 7505, 375),
    @(44295, 70, 8000, 8000, 8000, 400),
    @(44321, 250, 7000, 7000, 7000, 350),
    @(44412, 160, 7500, 8000, 7750, 388),
    @(44214, 50, 8000, 8000, 8000, 400),
    @(44314, 160, 8000, 8000, 8000, 400),
    @(44435, 302, 7000, 8000, 7500, 375),
    @(44405, 160, 7000, 8000, 7500, 375),
    @(44224, 120, 6000, 7000, 6667, 333),
    @(44391, 52, 7000, 8000, 7500, 375),
    @(44278, 70, 8000, 8000, 8000, 400),
    @(44358, 160, 7500, 8000, 7750, 388),
    @(44292, 70, 8000, 8000, 8000, 400),
    @(44419, 160, 7000, 8000, 7500, 375),
    @(44195, 70, 7000, 7000, 7000, 350),
    @(44265, 70, 8000, 8000, 8000, 400),
    @(44433, 142, 7000, 8000, 7500, 375),
    @(44344, 210, 8000, 8000, 8000, 400),
    @(44281, 250, 8000, 8000, 8000, 400),
    @(44160, 50, 7000, 8000, 7600, 380),
    @(44301, 160, 8000, 8000, 8000, 400),
    @(44272, 160, 8000, 8000, 8000, 400),
    @(44370, 160, 7500, 8000, 7750, 388),
    @(44232, 60, 7000, 7000, 7000, 350),
    @(44398, 70, 7500, 8000, 7750, 388),
    @(44286, 70, 8000, 8000, 8000, 400),
    @(44244, 70, 8000, 8000, 8000, 400),
    @(44208, 50, 8000, 8000, 8000, 400),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $firstRow + $i
    $vals = $data[$i]
    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("D$r").NumberFormat = $ws.Range("D2").NumberFormat
    $ws.Range("J$r").Value = $vals[1]
    $ws.Range("K$r").Value = $vals[2]
    $ws.Range("L$r").Value = $vals[3]
    $ws.Range("M$r").Value = $vals[4]
    $ws.Range("P$r").Value = $vals[5]
}
